$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update M2 (cited_by_count): 18 -> 19
$ws.Range("M2").Value = "'19"

# Row 5 and row 6 author/publication records were swapped (DOM and Banner author ids added)
$ws.Range("A5").Value = "'" + 'Danielle Campbell, Portia D. Cowlings, Martha Tholanah, Mallery Jenna Robinson, Gail E. Graham, Scovia Aseru, Karine Dubé, Susan E. Cohn, Katharine J. Bar, Elizabeth Connick, Rosie Mngqbisa, Eileen P. Scully, Jamila K. Stockman, Sara Gianella'
$ws.Range("A6").Value = "'" + 'Matthew T Ollerton, Joy M. Folkvord, Kristina K. Peachman, Soumya Shashikumar, Elaine Morrison, Linda L. Jagodzinski, Sheila A. Peel, Mohammad Khreiss, Richard T. D’Aquila, Sofía Casares, Mangala Rao, Elizabeth Connick'
$ws.Range("B5").Value = "'" + 'Charles R. Drew University of Medicine and Science (CDU), Los Angeles, CA; Joint Doctoral Program in Public Health, University of California, San Diego/San Diego State University, La Jolla, CA; Charles R. Drew University of Medicine and Science (CDU), Los Angeles, CA; Department of Education, Graduate School of Education and Psychology, Pepperdine University, Los Angeles, CA; Milton Park Clinical Research Site, Community Advisory Board, University of Zimbabwe Clinical Trials Research Center, Milton Park; ; Johns Hopkins University AIDS Clinical Trials Group Clinical Research Site, Community Advisory Board, Baltimore, MD; Joint Clinical Research Centre (JCRC)/Kampala Clinical Research Site, Kampala, Uganda; ; Infectious Diseases Division, Feinberg School of Medicine, Northwestern University, Chicago, IL; University of Pennsylvania, Philadelphia, PA; University of Arizona, Tuscon, AZ; ; Division of Infectious Diseases, Department of Medicine, Johns Hopkins University, Baltimore, MD; Division of Infectious Diseases and Global Public Health, Department of Medicine, University of California, San Diego (UCSD), La Jolla, CA; Division of Infectious Diseases and Global Public Health, Department of Medicine, University of California, San Diego (UCSD), La Jolla, CA'
$ws.Range("B6").Value = "'" + 'Department of Medicine, University of Arizona, Tucson, AZ, United States; Department of Medicine, University of Arizona, Tucson, AZ, United States; Laboratory of Adjuvant and Antigen Research, United States Military HIV Research Program, Walter Reed Army Institute of Research, Silver Spring, MD, United States; US Military Malaria Vaccine Program, Naval Medical Research Center, Silver Spring, MD, United States; Laboratory of Adjuvant and Antigen Research, United States Military HIV Research Program, Walter Reed Army Institute of Research, Silver Spring, MD, United States; Diagnostics and Countermeasure Branch, Walter Reed Army Institute of Research, Silver Spring, MD, United States; Diagnostics and Countermeasure Branch, Walter Reed Army Institute of Research, Silver Spring, MD, United States; Department of Surgery, University of Arizona, Tucson, AZ, United States; Department of Medicine, Feinberg School of Medicine, Northwestern University, Chicago, IL, United States; US Military Malaria Vaccine Program, Naval Medical Research Center, Silver Spring, MD, United States; Laboratory of Adjuvant and Antigen Research, United States Military HIV Research Program, Walter Reed Army Institute of Research, Silver Spring, MD, United States; Department of Medicine, University of Arizona, Tucson, AZ, United States'
$ws.Range("C5").Value = "'" + 'https://openalex.org/W4295048088'
$ws.Range("C6").Value = "'" + 'https://openalex.org/W4310040230'
$ws.Range("D5").Value = "'" + 'A Community Call to Action to Prioritize Inclusion and Enrollment of Women in HIV Cure-related Research'
$ws.Range("D6").Value = "'" + 'HIV-1 infected humanized DRAGA mice develop HIV-specific antibodies despite lack of canonical germinal centers in secondary lymphoid tissues'
$ws.Range("E5").Value = "'" + '2022-09-02'
$ws.Range("E6").Value = "'" + '2022-11-25'
$ws.Range("F5").Value = "'" + 'Journal of Acquired Immune Deficiency Syndromes'
$ws.Range("F6").Value = "'" + 'Frontiers in Immunology'
$ws.Range("G5").Value = "'" + 'Lippincott Williams & Wilkins'
$ws.Range("G6").Value = "'" + 'Frontiers Media'
$ws.Range("H5").Value = "'" + 'https://doi.org/10.1097/qai.0000000000003084'
$ws.Range("H6").Value = "'" + 'https://doi.org/10.3389/fimmu.2022.1047277'
$ws.Range("I5").Value = "'" + 'cc-by-nc-nd'
$ws.Range("I6").Value = "'" + 'cc-by'
$ws.Range("J5").Value = "'" + 'publishedVersion'
$ws.Range("J6").Value = "'" + 'publishedVersion'
$ws.Range("K5").Value = "'" + 'hybrid'
$ws.Range("K6").Value = "'" + 'gold'
$ws.Range("L5").Value = "'" + 'en'
$ws.Range("L6").Value = "'" + 'en'
$ws.Range("N5").Value = "'" + '2022'
$ws.Range("N6").Value = "'" + '2022'
$ws.Range("O5").Value = "'" + 'https://pubmed.ncbi.nlm.nih.gov/36083494'
$ws.Range("O6").Value = "'" + 'https://pubmed.ncbi.nlm.nih.gov/36505432'
$ws.Range("P5").Value = "'" + 'https://doi.org/10.1097/qai.0000000000003084'
$ws.Range("P6").Value = "'" + 'https://doi.org/10.3389/fimmu.2022.1047277'
$ws.Range("Q5").Value = "'" + 'article'
$ws.Range("Q6").Value = "'" + 'article'

# Update cited_by_count for rows 5 and 6 (M5: 2->3, M6: 1->2)
$ws.Range("M5").Value = "'3"
$ws.Range("M6").Value = "'2"